$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 changes from the text "R40" to the text "1" (still a text string,
# not a number -- the workbook keeps storing it as a shared string). Using a
# plain `.Value = "1"` assignment would make Excel infer a numeric type (or
# force a new "Text" number format onto the cell, changing its style id).
# Going through a TEXT() formula and then converting that formula to a
# static value via Copy / PasteSpecial(values) keeps the cell's existing
# style untouched while still leaving a text ("s") cell behind, matching
# the original cell's type.
$ws.Range("B11").Formula = "=TEXT(1,""0"")"
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

$wb.Save()
